$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-04 Friday" "2025-04-05 Saturday"

Replace-Text "518÷7=74, 0" "772÷6=128, 4"
Replace-Text "167÷8=20, 7" "416÷7=59, 3"
Replace-Text "255÷3=85, 0" "687÷7=98, 1"
Replace-Text "544÷5=108, 4" "105÷9=11, 6"
Replace-Text "759÷3=253, 0" "570÷9=63, 3"
Replace-Text "972÷4=243, 0" "212÷6=35, 2"
Replace-Text "630÷6=105, 0" "647÷2=323, 1"
Replace-Text "333÷7=47, 4" "560÷2=280, 0"
Replace-Text "285÷6=47, 3" "245÷3=81, 2"
Replace-Text "319÷9=35, 4" "930÷2=465, 0"
Replace-Text "249÷8=31, 1" "143÷3=47, 2"
Replace-Text "363÷5=72, 3" "735÷8=91, 7"
Replace-Text "672÷2=336, 0" "274÷9=30, 4"
Replace-Text "127÷8=15, 7" "991÷4=247, 3"
Replace-Text "295÷4=73, 3" "284÷3=94, 2"
Replace-Text "975÷7=139, 2" "222÷7=31, 5"
Replace-Text "885÷5=177, 0" "200÷6=33, 2"
Replace-Text "589÷7=84, 1" "596÷4=149, 0"
Replace-Text "442÷7=63, 1" "860÷7=122, 6"
Replace-Text "714÷2=357, 0" "971÷9=107, 8"
Replace-Text "315÷5=63, 0" "507÷4=126, 3"
Replace-Text "911÷4=227, 3" "705÷2=352, 1"
Replace-Text "974÷8=121, 6" "940÷6=156, 4"
Replace-Text "808÷5=161, 3" "133÷3=44, 1"
Replace-Text "905÷7=129, 2" "151÷2=75, 1"
